$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Values are stored as text in the sheet (e.g. '311.86', '-7.96%'), so we
# prefix with an apostrophe to force text entry (avoids numeric/percent
# auto-conversion), then reset the style back to Normal so no stray
# quote-prefix formatting is left behind on the cell.
$ws.Range("D2").Value = "'311.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'40.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-7.96%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.112"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.07853"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-5.86%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.330"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-2.23%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.672"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-13.98%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9244"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-4.78%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1073"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-4.89%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1776"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-5.66%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09108"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-6.05%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04437"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-3.71%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'7.182"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-17.82%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'-0.14%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001289"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.72%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005942"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.85%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.345"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.62%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.559"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.71%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.80%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1383"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.83%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'2.94%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04170"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.36%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001246"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.11%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004145"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-6.37%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001230"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-5.61%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003002"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.56%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02431"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-9.69%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05283"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-4.96%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.008015"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.97%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'-3.80%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007491"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'2.48%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002033"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.82%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008241"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'4.94%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3107"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-11.12%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006786"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.29%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000756"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.58%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003424"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-1.94%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004132"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'16.80%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002117"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.58%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002016"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.58%"
$ws.Range("E51").Style = "Normal"
